$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value "wewefwef" in cell E14 (this also introduces a new
# shared-string entry and extends the sheet's used range/dimension).
$ws.Range("E14").Value = "wewefwef"

# Move/update the active selection to the newly written cell, matching
# the author's final cursor position.
[void]$ws.Range("E14").Select()
